$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "65.076.55"
Set-TextValue $ws.Range("E2") "  +1.89%  "
Set-TextValue $ws.Range("D3") "3.180.83"
Set-TextValue $ws.Range("E3") "  +4.15%  "
Set-TextValue $ws.Range("E4") "  -0.06%  "
Set-TextValue $ws.Range("D5") "579.96"
Set-TextValue $ws.Range("E5") "  +3.67%  "
Set-TextValue $ws.Range("D6") "151.86"
Set-TextValue $ws.Range("E6") "  +6.43%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.02%  "
Set-TextValue $ws.Range("D8") "3.179.57"
Set-TextValue $ws.Range("E8") "  +4.18%  "
Set-TextValue $ws.Range("D9") "0.534"
Set-TextValue $ws.Range("E9") "  +3.82%  "
Set-TextValue $ws.Range("E10") "  +5.98%  "
Set-TextValue $ws.Range("D11") "6.24"
Set-TextValue $ws.Range("E11") "  -0.38%  "
Set-TextValue $ws.Range("E12") "  +2.40%  "
Set-TextValue $ws.Range("E13") "  +18.32%  "
Set-TextValue $ws.Range("D14") "37.95"
Set-TextValue $ws.Range("E14") "  +6.01%  "
Set-TextValue $ws.Range("D15") "3.698.30"
Set-TextValue $ws.Range("E15") "  +4.05%  "
Set-TextValue $ws.Range("D16") "65.123.35"
Set-TextValue $ws.Range("E16") "  +1.84%  "
Set-TextValue $ws.Range("D17") "3.175.84"
Set-TextValue $ws.Range("E17") "  +3.85%  "
Set-TextValue $ws.Range("D18") "7.19"
Set-TextValue $ws.Range("E18") "  +5.63%  "
Set-TextValue $ws.Range("E19") "  +1.58%  "
Set-TextValue $ws.Range("D20") "514.52"
Set-TextValue $ws.Range("E20") "  +7.98%  "
Set-TextValue $ws.Range("D21") "14.92"
Set-TextValue $ws.Range("E21") "  +5.78%  "
Set-TextValue $ws.Range("D22") "0.732"
Set-TextValue $ws.Range("E22") "  +7.02%  "
Set-TextValue $ws.Range("D23") "15.18"
Set-TextValue $ws.Range("E23") "  +3.58%  "
Set-TextValue $ws.Range("D24") "7.82"
Set-TextValue $ws.Range("E24") "  +3.60%  "
Set-TextValue $ws.Range("D25") "85.49"
Set-TextValue $ws.Range("E25") "  +3.24%  "
Set-TextValue $ws.Range("E26") "  -0.11%  "
Set-TextValue $ws.Range("D27") "9.01"
Set-TextValue $ws.Range("E27") "  +9.90%  "
Set-TextValue $ws.Range("D28") "2.94"
Set-TextValue $ws.Range("E28") "  +4.70%  "
Set-TextValue $ws.Range("E29") "  +7.47%  "
Set-TextValue $ws.Range("D30") "27.98"
Set-TextValue $ws.Range("E30") "  +6.40%  "
Set-TextValue $ws.Range("E31") "  +13.68%  "
Set-TextValue $ws.Range("E32") "  -0.03%  "
Set-TextValue $ws.Range("E33") "  +5.23%  "
Set-TextValue $ws.Range("D34") "6.33"
Set-TextValue $ws.Range("E34") "  +9.70%  "
Set-TextValue $ws.Range("D35") "6.61"
Set-TextValue $ws.Range("E35") "  +6.10%  "
Set-TextValue $ws.Range("D36") "55.71"
Set-TextValue $ws.Range("E36") "  +2.10%  "
Set-TextValue $ws.Range("D37") "0.0902"
Set-TextValue $ws.Range("E37") "  +10.64%  "
Set-TextValue $ws.Range("B38") "Bittensor"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D38") "476.05"
Set-TextValue $ws.Range("E38") "  +6.35%  "
Set-TextValue $ws.Range("B39") "dogwifhat"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D39") "3.15"
Set-TextValue $ws.Range("E39") "  +11.66%  "
Set-TextValue $ws.Range("E40") "  +3.12%  "
Set-TextValue $ws.Range("E41") "  +4.66%  "
Set-TextValue $ws.Range("D42") "3.077.82"
Set-TextValue $ws.Range("E42") "  +2.01%  "
Set-TextValue $ws.Range("E43") "  +2.46%  "
Set-TextValue $ws.Range("E44") "  +6.20%  "
Set-TextValue $ws.Range("E45") "  +5.61%  "
Set-TextValue $ws.Range("D46") "29.18"
Set-TextValue $ws.Range("E46") "  +3.02%  "
Set-TextValue $ws.Range("E47") "  +19.50%  "
Set-TextValue $ws.Range("E49") "  +2.32%  "
Set-TextValue $ws.Range("E50") "  +8.27%  "
Set-TextValue $ws.Range("D51") "120.59"
Set-TextValue $ws.Range("E51") "  +2.49%  "
